$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Push the old row 6 ("start_war" / "468, 440" trigger row) down to row 7 ---
# Copy formats first, then values, so the brand-new row 7 picks up the same
# cell styles (s="1"/"2") that row 6 used to have.
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Turn (the now free) row 6 into a new "apply_team" step ---
# Give it the same gray-filled / bordered look as the other option rows (2-5)
# by copying row 5's formatting before it gets overwritten below.
$ws.Range("A5:D5").Copy()
$ws.Range("A6:D6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A6").Value = "apply_team"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""

# --- Update the flow of rows 3-5 ---
$ws.Range("A3").Value = "team2"

$ws.Range("B5").Value = "airport"
$ws.Range("A5").Value = "468, 440"

$ws.Range("A4").Value = "apply_team"
$ws.Range("B4").Value = ""

# --- Update the active selection like the saved workbook shows ---
$ws.Range("A4").Select()
